$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 816, pushing existing rows 816:909 down to 818:911.
$ws.Rows.Item(816).Resize(2).Insert()

# Fill the two new rows (816-817) with the new weekly data point (date 44946).
$ws.Range("A816:A817").Value = 3
$ws.Range("B816:B817").Value = "Femacal de La Calera"
$ws.Range("C816:C817").Value = "Coquimbo"
$ws.Range("D816").Value = 44946
$ws.Range("D817").Value = 44946
$ws.Range("E816:E817").Value = 5
$ws.Range("F816:F817").Value = 100112008
$ws.Range("G816:G817").Value = "Coliflor"
$ws.Range("H816:H817").Value = "Sin especificar"
$ws.Range("I816").Value = "Primera"
$ws.Range("I817").Value = "Segunda"

$ws.Range("J816").Value = 3100
$ws.Range("K816").Value = 1000
$ws.Range("L816").Value = 1100
$ws.Range("M816").Value = 1052
$ws.Range("P816").Value = 1052

$ws.Range("J817").Value = 1400
$ws.Range("K817").Value = 900
$ws.Range("L817").Value = 900
$ws.Range("M817").Value = 900
$ws.Range("P817").Value = 900

$ws.Range("N816:N817").Value = "$/unidad"
$ws.Range("O816:O817").Value = "Provincia de Quillota"
$ws.Range("Q816:Q817").Value = 1
$ws.Range("R816:R817").Value = "Hortaliza"

# Match the date cell style used by the rest of column D.
$ws.Range("D816:D817").NumberFormat = $ws.Range("D818").NumberFormat
